$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values to repull/recalculate the data
$ws.Range("F2").Value = 2
$ws.Range("F3").Value = 2
$ws.Range("F4").Value = 3
$ws.Range("F6").Value = -1
$ws.Range("F9").Value = -1
